$d = $word.ActiveDocument

# --- locate donor ranges that already carry the exact character formatting
# we need to reproduce (this engine's Font.Name setter only ever emits
# w:ascii/w:hAnsi, never w:eastAsia/w:cs, so we clone fully-formed runs via
# FormattedText instead of rebuilding rPr from individual Font.* setters).
$donorCorr = $d.Content
$null = $donorCorr.Find.Execute("<cont/>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$donorAdd = $d.Content
$null = $donorAdd.Find.Execute("<div>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$donorPlainRtl = $d.Content
$null = $donorPlainRtl.Find.Execute("pinceau", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$ftCorr = $donorCorr.FormattedText
$ftAdd = $donorAdd.FormattedText
$ftPlainRtl = $donorPlainRtl.FormattedText

# --- find the unique anchor and split point: right before "est fort cler &"
# (which is itself immediately preceded by a single space in the existing
# run " est fort cler &"). Inserting here keeps that leading space, and the
# trailing "est fort cler &" text, as two pieces of the original run.
$anchor = $d.Content
$found = $anchor.Find.Execute("est fort cler &", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "anchor text 'est fort cler &' not found" }
$anchor.Collapse(1)
$pos = $anchor.Start

function Insert-PlainText($pos, $text) {
    $ins = $d.Range($pos, $pos)
    $ins.InsertBefore($text)
    $r = $d.Range($pos, $pos + $text.Length)
    $r.Font.Color = 0x000000
    return $pos + $text.Length
}

function Insert-Cloned($pos, $ft, $text) {
    $ins = $d.Range($pos, $pos)
    $ins.FormattedText = $ft
    $r = $d.Range($pos, $pos + $text.Length)
    $r.Text = $text
    return $pos + $text.Length
}

$pos = Insert-Cloned $pos $ftCorr "<corr>"
$pos = Insert-PlainText $pos "e"
$pos = Insert-Cloned $pos $ftCorr "</corr>"
$pos = Insert-Cloned $pos $ftAdd "<add>"
$pos = Insert-PlainText $pos "t"
$pos = Insert-Cloned $pos $ftAdd "</add>"
$pos = Insert-Cloned $pos $ftPlainRtl " "
